# "change monster to entity"
#
# Sheet1's column G held the "(S)MonsterFile(S)" / monster spawn-config
# fields; column H held the separate "(S)NPCFile(S)" / npc config fields.
# Monster + NPC configs were merged into a single generic "entity" config,
# so column G is renamed/retargeted to the entity wording and column H's
# now-obsolete content is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "(S)EntityFile(S)"
$ws.Range("G2").Value = "非玩家配表"
$ws.Range("G3").Value = "scripts\data\spawnpoints\entity_newplayermapinfo.xml"

$ws.Range("H1:H3").ClearContents()
